# Updates cryptos list prices/volumes, and reorders the Aave/Cronos/FraxShare rows
# to match the latest scrape (commit: "Updated cryptos list ... with GitHub Actions").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Worksheet, [string]$Addr, [string]$Val)
    $cell = $Worksheet.Range($Addr)
    # Force the cell to stay a text value (avoids Excel auto-converting
    # numeric-looking strings like "257.22" into a float), then strip the
    # temporary text number-format so no extra style is left behind.
    $cell.NumberFormat = "@"
    $cell.Value2 = $Val
    $cell.ClearFormats()
}

$updates = @{
    2  = @{ D = "43.616.24";  E = "  +2.33%  " }
    3  = @{ D = "2.197.80";   E = "  -0.26%  " }
    4  = @{ E = "  +0.11%  " }
    5  = @{ D = "257.22";     E = "  +1.19%  " }
    6  = @{ D = "84.44";      E = "  +12.89%  " }
    7  = @{ D = "0.620";      E = "  +0.86%  " }
    8  = @{ E = "  +0.06%  " }
    9  = @{ D = "0.593";      E = "  +0.89%  " }
    10 = @{ D = "44.83";      E = "  +9.06%  " }
    11 = @{ D = "0.0917";     E = "  +0.76%  " }
    12 = @{ D = "7.23";       E = "  +5.95%  " }
    13 = @{ E = "  +2.62%  " }
    14 = @{ D = "2.530.80";   E = "  -0.19%  " }
    15 = @{ D = "14.36";      E = "  +0.73%  " }
    16 = @{ D = "2.231.41";   E = "  +1.78%  " }
    17 = @{ D = "0.783";      E = "  +0.84%  " }
    18 = @{ D = "43.572.71";  E = "  +2.41%  " }
    19 = @{ E = "  +0.73%  " }
    20 = @{ D = "69.75";      E = "  -1.68%  " }
    21 = @{ D = "5.90";       E = "  -0.33%  " }
    22 = @{ D = "2.37";       E = "  +8.49%  " }
    23 = @{ D = "231.13";     E = "  +1.91%  " }
    24 = @{ D = "8.92";       E = "  -4.72%  " }
    26 = @{ D = "10.62";      E = "  +1.01%  " }
    27 = @{ E = "  +3.58%  " }
    28 = @{ D = "2.27";       E = "  +3.20%  " }
    29 = @{ D = "38.92";      E = "  -0.43%  " }
    30 = @{ E = "  +2.19%  " }
    31 = @{ D = "173.32";     E = "  +0.14%  " }
    32 = @{ D = "20.35";      E = "  +1.02%  " }
    33 = @{ D = "0.0861";     E = "  +2.87%  " }
    34 = @{ D = "5.30";       E = "  +2.22%  " }
    35 = @{ E = "  +1.58%  " }
    36 = @{ D = "0.110";      E = "  +1.39%  " }
    37 = @{ D = "0.0358";     E = "  +4.37%  " }
    38 = @{ D = "4.47";       E = "  +4.42%  " }
    39 = @{ D = "12.44";      E = "  +0.22%  " }
    40 = @{ D = "2.86";       E = "  +4.75%  " }
    41 = @{ E = "  +0.06%  " }
    42 = @{ D = "63.11";      E = "  +5.80%  " }
    43 = @{ D = "5.46";       E = "  +4.25%  " }
    44 = @{ D = "0.198";      E = "  +0.89%  " }
    45 = @{ B = "FraxShare"; C = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"; D = "8.32";    E = "  -0.38%  " }
    46 = @{ B = "Aave";      C = "https://coinranking.com/coin/ixgUfzmLR+aave-aave";    D = "100.07";  E = "  -1.13%  " }
    47 = @{ B = "Cronos";    C = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"; D = "0.0978"; E = "  +0.25%  " }
    48 = @{ E = "  +4.37%  " }
    49 = @{ E = "  +0.80%  " }
    50 = @{ D = "0.435";      E = "  -5.46%  " }
    51 = @{ D = "1.48";       E = "  +4.31%  " }
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    foreach ($col in $vals.Keys) {
        Set-TextValue $ws "$col$row" $vals[$col]
    }
}
